$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "q1"
$ws.Range("D1").Value = "q2"
$ws.Range("E1").Value = "q3"
$ws.Range("F1").Value = "p1"
$ws.Range("G1").Value = "p2"
$ws.Range("H1").Value = "p3"
$ws.Range("I1").Value = "tf"

$ws.Range("I2").Select()
